$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert three new columns before column H (old H,I,J -> new K,L,M)
$ws.Range("H1:J1").EntireColumn.Insert()

# Copy the header style (bold/centered/bordered) used by existing headers onto the
# newly inserted header cells.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 2. Rewrite the header row with the new (lower-cased / renamed) labels and the
#    three newly inserted yield columns.
$ws.Range("A1").Value = "name of instrument"
$ws.Range("B1").Value = "isin"
$ws.Range("C1").Value = "coupon"
$ws.Range("D1").Value = "industry"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "market value (mkt)"
$ws.Range("G1").Value = "% to net assets (nav)"
$ws.Range("H1").Value = "yield"
$ws.Range("I1").Value = "yield to call (ytc)"
$ws.Range("J1").Value = "yield to maturity (ytm)"
$ws.Range("K1").Value = "Type"
$ws.Range("L1").Value = "Scheme"
$ws.Range("M1").Value = "AmcName"

# 3. Move the "Sovereign" rating value from column C to column D (industry) for
#    each data row, leaving column C (coupon) blank.
$ws.Range("D2").Value = $ws.Range("C2").Value2
$ws.Range("C2").Value = ""
$ws.Range("D3").Value = $ws.Range("C3").Value2
$ws.Range("C3").Value = ""

# 4. Clear the stray values left behind in the newly inserted H/I columns, fill in
#    the new "yield to maturity (ytm)" values, and rewrite the "Type" column text
#    that landed in column K after the insert.
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "'0.067986"
$ws.Range("K2").Value = "Debt Instruments  NAN nan nan nan nan nan nan nan"

$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "'0.067983"
$ws.Range("K3").Value = "Debt Instruments  NAN nan nan nan nan nan nan nan"
